$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (34) down onto the
# three new rows so every cell picks up the exact same style indices
# (font/numberformat) already used for dates, datetimes and plain values.
$ws.Range("A34:I34").Copy()
$ws.Range("A35:I37").PasteSpecial(-4122)

# --- Row 35: 2019-03-28 (Thu) ---
$ws.Cells.Item(35, 1).Value = 43552.0
$ws.Cells.Item(35, 2).Value = "Thu"
$ws.Cells.Item(35, 3).Value = 2.8
$ws.Cells.Item(35, 4).Value = 3.0
$ws.Cells.Item(35, 5).Value = 43552.65277777778
$ws.Cells.Item(35, 6).Value = 43553.125
$ws.Cells.Item(35, 7).Value = 1000.0
$ws.Cells.Item(35, 8).Value = 0.0
$ws.Cells.Item(35, 9).Value = 342.0

# --- Row 36: 2019-04-26 (Fri) ---
$ws.Cells.Item(36, 1).Value = 43581.0
$ws.Cells.Item(36, 2).Value = "Fri"
$ws.Cells.Item(36, 3).Value = 3.0
$ws.Cells.Item(36, 4).Value = 3.0
$ws.Cells.Item(36, 5).Value = 43581.75
$ws.Cells.Item(36, 6).Value = 43582.28472222222
$ws.Cells.Item(36, 7).Value = 520.0
$ws.Cells.Item(36, 8).Value = 0.0
$ws.Cells.Item(36, 9).Value = 274.0

# --- Row 37: 2019-05-01 (Wed) ---
$ws.Cells.Item(37, 1).Value = 43586.0
$ws.Cells.Item(37, 2).Value = "Wed"
$ws.Cells.Item(37, 3).Value = 2.0
$ws.Cells.Item(37, 4).Value = 2.0
$ws.Cells.Item(37, 5).Value = 43586.575
$ws.Cells.Item(37, 6).Value = 43586.947916666664
$ws.Cells.Item(37, 7).Value = 260.0
$ws.Cells.Item(37, 8).Value = 0.0
$ws.Cells.Item(37, 9).Value = 1233.0

# The May-1st row (37) gets its own date format: a new number format
# "mmmm d" (full month name) distinct from the "mmm d" used elsewhere.
$ws.Range("A37").NumberFormat = "mmmm d"

Write-Host "rows 35-37 added"
